$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff represents a cyclic rotation of the weekly price records held in
# rows 2, 3, 5, 6 and 8 (columns D, M, N, O, P, S). Row 4 and 7 are untouched.
# Capture the "before" values first, then write them back in rotated order so
# we don't clobber data we still need to read.

$rows = 2, 3, 5, 6, 8
$cols = "D", "M", "N", "O", "P", "S"

$original = @{}
foreach ($r in $rows) {
    $original[$r] = @{}
    foreach ($c in $cols) {
        $original[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# new_row[target] = old_row[source]
$rotation = @{ 5 = 2; 6 = 5; 3 = 6; 8 = 3; 2 = 8 }

foreach ($target in $rotation.Keys) {
    $source = $rotation[$target]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $original[$source][$c]
    }
}
